# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose underlying data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    5  = 4
    15 = -3
    16 = -7
    17 = 6
    19 = 3
    21 = -6
    24 = -4
    25 = -1
    26 = -4
    29 = -3
    30 = 13
    32 = 9
    39 = -2
    40 = 4
    41 = -4
    42 = -3
    46 = -7
    48 = 4
    49 = -3
    52 = 5
    54 = -1
    55 = 0
    56 = -2
    57 = 5
    58 = -6
    60 = 4
    64 = -9
    70 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
